# storage host presentation final refactoring
#
# This script reproduces, via Excel COM automation, the changes captured in the
# target OOXML diff:
#   1. Workbook-level window view tweaks (window size, active tab).
#   2. "report" sheet: move the "value" header label from C1 to E1 (with its
#      formatting), make "report" the selected/active tab, and move the
#      selection to E10.
#   3. "service_tables" sheet: de-select it as the active tab, move the
#      selection to F97:H98, and zero out the G/H ("is active"/"is visible")
#      flag columns for a number of rows.

$wb = $excel.ActiveWorkbook

$wsReport = $wb.Worksheets.Item("report")
$wsServiceTables = $wb.Worksheets.Item("service_tables")

# ---------------------------------------------------------------------------
# 1. "report" sheet (Sheet1): relocate the "value" label from C1 to E1,
#    keeping its original formatting, then select it as the active sheet
#    with E10 as the active cell.
# ---------------------------------------------------------------------------
$wsReport.Activate()

$wsReport.Range("C1").Copy($wsReport.Range("E1"))
$wsReport.Range("C1").Clear()

$wsReport.Range("E10").Select()

# ---------------------------------------------------------------------------
# 2. "service_tables" sheet (Sheet2): zero out the G/H flag columns for the
#    affected rows, then move the selection to F97:H98 (leaving the existing
#    frozen header pane untouched).
# ---------------------------------------------------------------------------
$wsServiceTables.Activate()

$rowsGAndH = @(47,50,51,52,53,54,57,58,59,68,69,70,71,72,73,74,75,76,78,79,80,87,88,89,90,91,92,93,94,95,100,102,105)
$rowsHOnly = @(56,60,61,83,84,85,86,101,104)

foreach ($r in $rowsGAndH) {
    $wsServiceTables.Range("G$r").Value = 0
    $wsServiceTables.Range("H$r").Value = 0
}
foreach ($r in $rowsHOnly) {
    $wsServiceTables.Range("H$r").Value = 0
}

$wsServiceTables.Range("F97:H98").Select()

# ---------------------------------------------------------------------------
# 3. Re-activate "report" so it becomes the workbook's selected/active tab,
#    matching the final workbook view state.
# ---------------------------------------------------------------------------
$wsReport.Activate()
